$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Week 8 (row 16) actuals: Direct Hours and Team Hours, with their running
# cumulative totals computed the same way as the earlier weeks (rows 14-15).
$ws.Range("G16").Value = 15.5
$ws.Range("I16").Formula = "=G16+I15"
$ws.Range("J16").Value = 50.3
$ws.Range("L16").Formula = "=J16+L15"

# Leave the selection where Excel would land after entering the last value.
[void]$ws.Range("M16").Select()
